# BrowserStack Firefox browser code added
# Updates the TestData.xlsx "productName" test-data rows:
#   - productName_1 (row 7)  : ZARA COAT 3   -> ADIDAS ORIGINAL
#   - productName_5 (row 15) : ZARA COAT 3   -> ADIDAS ORIGINAL
#   - productName_6 (row 19) : ADIDAS ORIGINAL -> IPHONE 13 PRO
# and leaves the cursor/selection on the newly edited cell (B19),
# matching the saved workbook view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "ADIDAS ORIGINAL"
$ws.Range("B15").Value = "ADIDAS ORIGINAL"
$ws.Range("B19").Value = "IPHONE 13 PRO"

# B3/B4 picked up an explicit (but visually no-op) fill flag at some point;
# re-harmonize their formatting with B2, which already carries the clean
# "no fill" version of the same font/alignment combo.
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B19").Select()
